$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowByLabel($label) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value() -eq $label) {
            return $r
        }
    }
    return -1
}

# --- Remove whole rows that were dropped from the dataset ---
$rowsToDelete = @("RM 232", "SC 92")
foreach ($label in $rowsToDelete) {
    $r = Get-RowByLabel $label
    if ($r -gt 0) {
        $ws.Rows.Item($r).Delete()
    }
}

# --- Apply per-cell value changes (numbers set/updated, others cleared to missing) ---
$cellEdits = @{
    "RM 14"  = @{ "F" = $null }
    "RM 21"  = @{ "E" = -5.7 }
    "RM 38"  = @{ "E" = $null }
    "RM 58"  = @{ "F" = 17.65 }
    "RM 125" = @{ "E" = -6.5; "F" = $null }
    "RM 135" = @{ "E" = $null }
    "RM 140" = @{ "E" = -7; "F" = 16.48 }
    "RM 145" = @{ "F" = 16.6 }
    "SC 5"   = @{ "B" = $null }
    "SC 101" = @{ "B" = -20.4; "E" = $null; "F" = $null }
    "SC 119" = @{ "B" = $null; "E" = -6.8; "F" = $null }
    "SC 120" = @{ "F" = 16.89 }
    "SC 232" = @{ "F" = 17.53 }
}

foreach ($label in $cellEdits.Keys) {
    $r = Get-RowByLabel $label
    if ($r -gt 0) {
        $colEdits = $cellEdits[$label]
        foreach ($col in $colEdits.Keys) {
            $val = $colEdits[$col]
            $cellRef = "$col$r"
            if ($null -eq $val) {
                $ws.Range($cellRef).ClearContents()
            } else {
                $ws.Range($cellRef).Value = $val
            }
        }
    }
}
